$wb = $excel.ActiveWorkbook

# --- users sheet: add "icon" default value to column F for existing rows ---
$usersWs = $wb.Worksheets.Item("users")
$usersWs.Range("F2:F101").Value = "default_icon.png"

# Reflect the new selection / navigation state captured in the saved file:
# the users sheet becomes the active tab, scrolled down one row, with I11 selected.
[void]$usersWs.Activate()
$excel.ActiveWindow.ScrollRow = 2
[void]$usersWs.Range("I11").Select()

Write-Output "done"
